$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.298.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.46%  "
$ws.Range("E2").ClearFormats()

$ws.Range("D3").Value = "'3.427.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +3.00%  "
$ws.Range("E3").ClearFormats()

$ws.Range("E4").Value = "'  +0.37%  "
$ws.Range("E4").ClearFormats()

$ws.Range("D5").Value = "'547.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +3.15%  "
$ws.Range("E5").ClearFormats()

$ws.Range("D6").Value = "'178.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.97%  "
$ws.Range("E6").ClearFormats()

$ws.Range("D7").Value = "'0.634"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +7.01%  "
$ws.Range("E7").ClearFormats()

$ws.Range("E8").Value = "'  +0.13%  "
$ws.Range("E8").ClearFormats()

$ws.Range("D9").Value = "'0.621"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.75%  "
$ws.Range("E9").ClearFormats()

$ws.Range("D10").Value = "'0.150"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +7.75%  "
$ws.Range("E10").ClearFormats()

$ws.Range("D11").Value = "'53.17"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -1.66%  "
$ws.Range("E11").ClearFormats()

$ws.Range("D12").Value = "'0.0000268"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +2.96%  "
$ws.Range("E12").ClearFormats()

$ws.Range("D13").Value = "'9.10"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.92%  "
$ws.Range("E13").ClearFormats()

$ws.Range("D14").Value = "'3.984.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +3.80%  "
$ws.Range("E14").ClearFormats()

$ws.Range("E15").Value = "'  +2.18%  "
$ws.Range("E15").ClearFormats()

$ws.Range("D16").Value = "'3.434.10"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +3.69%  "
$ws.Range("E16").ClearFormats()

$ws.Range("D17").Value = "'18.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +3.35%  "
$ws.Range("E17").ClearFormats()

$ws.Range("D18").Value = "'65.338.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.38%  "
$ws.Range("E18").ClearFormats()

$ws.Range("D19").Value = "'11.72"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +3.84%  "
$ws.Range("E19").ClearFormats()

$ws.Range("D20").Value = "'0.975"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +1.30%  "
$ws.Range("E20").ClearFormats()

$ws.Range("D21").Value = "'412.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +7.21%  "
$ws.Range("E21").ClearFormats()

$ws.Range("D22").Value = "'3.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +6.92%  "
$ws.Range("E22").ClearFormats()

$ws.Range("D23").Value = "'4.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +1.99%  "
$ws.Range("E23").ClearFormats()

$ws.Range("D24").Value = "'84.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +2.66%  "
$ws.Range("E24").ClearFormats()

$ws.Range("D25").Value = "'10.72"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -4.45%  "
$ws.Range("E25").ClearFormats()

$ws.Range("D26").Value = "'2.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +3.10%  "
$ws.Range("E26").ClearFormats()

$ws.Range("B27").Value = "'LEO"
$ws.Range("B27").ClearFormats()
$ws.Range("C27").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C27").ClearFormats()
$ws.Range("D27").Value = "'6.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.20%  "
$ws.Range("E27").ClearFormats()

$ws.Range("B28").Value = "'InternetComputer(DFINITY)"
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "'12.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +6.15%  "
$ws.Range("E28").ClearFormats()

$ws.Range("D29").Value = "'8.83"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +6.59%  "
$ws.Range("E29").ClearFormats()

$ws.Range("D30").Value = "'29.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.71%  "
$ws.Range("E30").ClearFormats()

$ws.Range("B31").Value = "'Bittensor"
$ws.Range("B31").ClearFormats()
$ws.Range("C31").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C31").ClearFormats()
$ws.Range("D31").Value = "'610.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -5.01%  "
$ws.Range("E31").ClearFormats()

$ws.Range("B32").Value = "'NEARProtocol"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = "'6.44"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -4.74%  "
$ws.Range("E32").ClearFormats()

$ws.Range("D33").Value = "'11.56"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +2.65%  "
$ws.Range("E33").ClearFormats()

$ws.Range("D34").Value = "'0.108"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +2.01%  "
$ws.Range("E34").ClearFormats()

$ws.Range("D35").Value = "'58.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +2.74%  "
$ws.Range("E35").ClearFormats()

$ws.Range("D36").Value = "'0.147"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +16.13%  "
$ws.Range("E36").ClearFormats()

$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -0.02%  "
$ws.Range("E37").ClearFormats()

$ws.Range("D38").Value = "'36.98"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.93%  "
$ws.Range("E38").ClearFormats()

$ws.Range("D39").Value = "'0.0₃0772"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +1.75%  "
$ws.Range("E39").ClearFormats()

$ws.Range("D40").Value = "'0.375"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -2.43%  "
$ws.Range("E40").ClearFormats()

$ws.Range("D41").Value = "'3.132.54"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +4.91%  "
$ws.Range("E41").ClearFormats()

$ws.Range("D42").Value = "'3.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +2.19%  "
$ws.Range("E42").ClearFormats()

$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.44%  "
$ws.Range("E43").ClearFormats()

$ws.Range("D44").Value = "'2.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -4.19%  "
$ws.Range("E44").ClearFormats()

$ws.Range("D45").Value = "'2.78"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +3.43%  "
$ws.Range("E45").ClearFormats()

$ws.Range("B46").Value = "'WEMIXToken"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'2.72"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +1.24%  "
$ws.Range("E46").ClearFormats()

$ws.Range("B47").Value = "'VeChain"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'0.0407"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.84%  "
$ws.Range("E47").ClearFormats()

$ws.Range("B48").Value = "'ApeXProtocol"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'3.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +2.01%  "
$ws.Range("E48").ClearFormats()

$ws.Range("D49").Value = "'0.130"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +3.46%  "
$ws.Range("E49").ClearFormats()

$ws.Range("D50").Value = "'138.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +0.19%  "
$ws.Range("E50").ClearFormats()

$ws.Range("D51").Value = "'8.33"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +0.33%  "
$ws.Range("E51").ClearFormats()

